$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per row, resulting from a data repull.
$updates = @{
    2  = -1
    3  = 2
    4  = -1
    5  = 5
    6  = -2
    7  = 8
    8  = 2
    13 = -2
    14 = 3
    15 = -3
    16 = 10
    17 = 3
    18 = -1
    22 = 11
    23 = 1
    24 = 3
    25 = -1
    27 = 5
    28 = -1
    29 = -1
    30 = 8
    31 = 0
    33 = -1
    34 = 7
    35 = 2
    37 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
